$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (cohort 2019, period 6): num_customers 8 -> 9, retention_rate recalculated
$ws.Range("C16").Value = 9
$ws.Range("E16").Value = 0.004271476032273375

# Row 27 (cohort 2022, period 4): num_customers 39 -> 41, retention_rate recalculated
$ws.Range("C27").Value = 41
$ws.Range("E27").Value = 0.01820603907637655

# Row 31 (cohort 2023, period 3): num_customers 40 -> 41, retention_rate recalculated
$ws.Range("C31").Value = 41
$ws.Range("E31").Value = 0.01773356401384083

# Row 36 (cohort 2024, period 1): num_customers 111 -> 112, retention_rate recalculated
$ws.Range("C36").Value = 112
$ws.Range("E36").Value = 0.05803108808290156

# Row 37 (cohort 2025, period 0): num_customers/cohort_size 682 -> 692, retention_rate stays 1
$ws.Range("C37").Value = 692
$ws.Range("D37").Value = 692
$ws.Range("E37").Value = 1
